$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 16 from serial date
# 45185 (2023-09-16) to 45204 (2023-10-05), matching the diff.
for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45185) {
        $cell.Value = 45204
    }
}
